$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range (A1:D3) so stale cells (C,D columns, row's beyond 5) are removed
$ws.Range("A1:D3").Clear()

# New header + data for underwater/weather flag codes
$data = @(
    @("flag", "description"),
    @("w", "do sensor cleaned"),
    @("n", "weather station adjusted to orient more northward"),
    @("R", "removed buoy"),
    @("D", "buoy deployed for season")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Column B width adjustment (no longer auto bestFit, fixed width 19.5 in OOXML units)
$ws.Columns.Item(2).ColumnWidth = 18.86

# Update the active selection to match the target state
$ws.Range("D11").Select()
